# Applies the "script update" diff to the Copa de la Liga Profesional 2023
# sheet:
#   - four pairs of adjacent rows had their match data (columns F..V) swapped
#     (the "home"/"away" pairing in the source feed got re-ordered, while the
#     row's Indice/pais/torneio/temporada/data_partida in columns A..E stayed
#     put), and
#   - one brand-new match (row 155) was appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F..V hold the actual match data (home team .. url); A..E (Indice,
# pais, torneio, temporada, data_partida) are left alone.
$matchCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-MatchRows {
    param([int]$row1, [int]$row2)

    $left = @{}
    $right = @{}
    foreach ($col in $matchCols) {
        $left[$col] = $ws.Range("$col$row1").Value2
        $right[$col] = $ws.Range("$col$row2").Value2
    }
    foreach ($col in $matchCols) {
        $ws.Range("$col$row1").Value2 = $right[$col]
        $ws.Range("$col$row2").Value2 = $left[$col]
    }
}

Swap-MatchRows 98 99
Swap-MatchRows 106 107
Swap-MatchRows 125 126
Swap-MatchRows 136 137

# Append the new match as row 155, mirroring the formatting of the previous
# last row (154) so the bold/boxed Indice cell and the datetime-formatted
# data_partida cell keep their existing styles instead of minting new ones.
$ws.Range("A154:V154").Copy()
$ws.Range("A155:V155").PasteSpecial(-4122)

$ws.Range("A155").Value2 = 154
$ws.Range("B155").Value2 = "argentina"
$ws.Range("C155").Value2 = "copa-de-la-liga-profesional"
# "temporada" is stored as text ("2023") throughout the sheet, not a number;
# the leading apostrophe forces the COM layer to keep it a string instead of
# auto-coercing the numeric-looking literal.
$ws.Range("D155").Value2 = "'2023"
$ws.Range("E155").Value2 = 45232
$ws.Range("F155").Value2 = "Defensa y Justicia"
$ws.Range("G155").Value2 = 2
$ws.Range("H155").Value2 = "Racing Club"
$ws.Range("I155").Value2 = 2
$ws.Range("J155").Value2 = 2.92
$ws.Range("K155").Value2 = "26/10/2023 02:12"
$ws.Range("L155").Value2 = 2.88
$ws.Range("M155").Value2 = "01/11/2023 23:58"
$ws.Range("N155").Value2 = 3.32
$ws.Range("O155").Value2 = "26/10/2023 02:12"
$ws.Range("P155").Value2 = 3.45
$ws.Range("Q155").Value2 = "01/11/2023 23:51"
$ws.Range("R155").Value2 = 2.51
$ws.Range("S155").Value2 = "26/10/2023 02:12"
$ws.Range("T155").Value2 = 2.52
$ws.Range("U155").Value2 = "01/11/2023 23:58"
$ws.Range("V155").Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/defensa-y-justicia-racing-club/86NfV6b3/"
